$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "95.685.97"
$ws.Range("E2").Value = "  +1.06%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.566.74"
$ws.Range("E3").Value = "  +1.22%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.97"
$ws.Range("E5").Value = "  -1.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "652.55"
$ws.Range("E6").Value = "  +3.32%  "

$ws.Range("E7").Value = "  +1.85%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.401"
$ws.Range("E8").Value = "  +1.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("E10").Value = "  -0.33%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.564.45"
$ws.Range("E11").Value = "  +1.33%  "

$ws.Range("E12").Value = "  +0.41%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.28"
$ws.Range("E13").Value = "  -2.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.57"
$ws.Range("E14").Value = "  +4.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.228.84"
$ws.Range("E15").Value = "  +1.00%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.549.02"
$ws.Range("E16").Value = "  +1.06%  "

$ws.Range("E17").Value = "  +0.62%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.565.63"
$ws.Range("E18").Value = "  +1.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.92"
$ws.Range("E19").Value = "  -4.53%  "

$ws.Range("E20").Value = "  -0.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.88"
$ws.Range("E21").Value = "  -0.98%  "

$ws.Range("E22").Value = "  +3.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "507.49"
$ws.Range("E23").Value = "  -2.05%  "

$ws.Range("E24").Value = "  -4.17%  "

$ws.Range("E25").Value = "  +3.89%  "

$ws.Range("E26").Value = "  -1.66%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "92.16"
$ws.Range("E27").Value = "  -0.62%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.55"
$ws.Range("E28").Value = "  +2.42%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.756.57"
$ws.Range("E29").Value = "  +1.36%  "

$ws.Range("E30").Value = "  +3.78%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.32"
$ws.Range("E31").Value = "  -1.87%  "

$ws.Range("E32").Value = "  -0.09%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.140"
$ws.Range("E33").Value = "  -1.60%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.11%  "

$ws.Range("E35").Value = "  -2.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.97"
$ws.Range("E36").Value = "  +6.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.560"
$ws.Range("E37").Value = "  -0.31%  "

$ws.Range("E38").Value = "  +8.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "561.01"
$ws.Range("E39").Value = "  -4.38%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.49"
$ws.Range("E40").Value = "  +2.59%  "

$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("E42").Value = "  +0.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.905"
$ws.Range("E43").Value = "  -2.45%  "

$ws.Range("E44").Value = "  +3.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "35.11"
$ws.Range("E45").Value = "  +36.60%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.31"
$ws.Range("E46").Value = "  +6.59%  "

$ws.Range("E47").Value = "  +2.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.58"
$ws.Range("E48").Value = "  -0.77%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0414"
$ws.Range("E49").Value = "  -2.17%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.57"
$ws.Range("E50").Value = "  +0.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.47"
$ws.Range("E51").Value = "  -0.99%  "
